$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K (strikeout) values for column G, rows 2-33, regenerated from source data
$newValues = @{
    2  = 6
    3  = 4
    4  = 6
    5  = 9
    6  = 6
    7  = 6
    8  = 2
    9  = 6
    10 = 5
    11 = 10
    12 = 6
    13 = 8
    14 = 10
    15 = 10
    16 = 4
    17 = 10
    18 = 5
    19 = 4
    20 = 4
    21 = 9
    22 = 9
    23 = 3
    24 = 8
    25 = 8
    26 = 1
    27 = 3
    28 = 10
    29 = 9
    30 = 5
    31 = 3
    32 = 4
    33 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
